$d = $word.ActiveDocument

# Merge the split runs in the Title paragraph into a single run by doing a
# find/replace over the full text (this collapses the multiple w:r elements
# Word originally split the text into back into one run with the same text).
$d.Content.Find.Execute("Questions: Introduction to Matrices", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Questions: Introduction to Matrices", 2)

# Merge the split runs in the Author paragraph.
$d.Content.Find.Execute("Jessica Taberner", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Jessica Taberner", 2)

# Merge the split runs in the Abstract paragraph.
$d.Content.Find.Execute("A selection of questions on matrices.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "A selection of questions on matrices.", 2)
